# Fruta / hortaliza, semanal
# Insert a new weekly record at the top of the price history table
# (row 44), pushing the existing rows 44-73 down to 45-74.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 44:73 down by one row to make room for the new record.
$ws.Rows("44:44").Insert()

# Populate the newly inserted row 44 with this week's data.
$ws.Range("A44").Value = 10
$ws.Range("B44").Value = "Vega Modelo de Temuco"
$ws.Range("C44").Value = "La Araucanía"
$ws.Range("D44").Value = 44729
$ws.Range("E44").Value = 9
$ws.Range("F44").Value = "Fruta"
$ws.Range("G44").Value = 100108
$ws.Range("H44").Value = "Tropicales y subtropicales"
$ws.Range("I44").Value = 100108007
$ws.Range("J44").Value = "Coco"
$ws.Range("K44").Value = "Sin especificar"
$ws.Range("L44").Value = "Primera"
$ws.Range("M44").Value = 20
$ws.Range("N44").Value = 28000
$ws.Range("O44").Value = 28000
$ws.Range("P44").Value = 28000
$ws.Range("Q44").Value = "$/malla 20 unidades"
$ws.Range("R44").Value = "Perú"
$ws.Range("S44").Value = 1400
$ws.Range("T44").Value = 20
